$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "23.911.22", "1.003")
# are preserved exactly as text instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.911.22"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.646.89"
$ws.Range("E3").Value = "  +1.60%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "308.81"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "0.3888"
$ws.Range("E7").Value = "  -1.04%  "

$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  -0.76%  "

$ws.Range("D9").Value = "51.22"
$ws.Range("E9").Value = "  +3.25%  "

$ws.Range("D10").Value = "1.345"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").Value = "0.08419"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "23.78"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").Value = "7.058"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("D15").Value = "7.842"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("D16").Value = "0.00001312"
$ws.Range("E16").Value = "  +2.03%  "

$ws.Range("D17").Value = "1.652.60"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("D18").Value = "94.29"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").Value = "0.06984"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "19.64"
$ws.Range("E20").Value = "  -1.71%  "

$ws.Range("D21").Value = "6.912"
$ws.Range("E21").Value = "  +1.19%  "

$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "13.69"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").Value = "23.916.06"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "2.457"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").Value = "2.960"
$ws.Range("E26").Value = "  +4.46%  "

$ws.Range("D27").Value = "22.01"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").Value = "150.66"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("D29").Value = "5.427"
$ws.Range("E29").Value = "  +2.34%  "

$ws.Range("D30").Value = "138.24"
$ws.Range("E30").Value = "  -1.58%  "

$ws.Range("D31").Value = "7.760"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").Value = "2.499"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").Value = "1.834.02"
$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("D34").Value = "1.042"
$ws.Range("E34").Value = "  +5.18%  "

$ws.Range("D35").Value = "0.08049"
$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("D36").Value = "0.02951"
$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("D37").Value = "6.714"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("D38").Value = "10.85"
$ws.Range("E38").Value = "  +5.01%  "

$ws.Range("D39").Value = "0.2673"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").Value = "0.09109"
$ws.Range("E40").Value = "  -0.60%  "

$ws.Range("D41").Value = "0.7525"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").Value = "13.40"
$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").Value = "1.424"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").Value = "16.25"
$ws.Range("E44").Value = "  +1.52%  "

$ws.Range("D45").Value = "0.6921"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").Value = "2.449"
$ws.Range("E46").Value = "  -1.05%  "

$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").Value = "0.08267"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").Value = "134.11"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").Value = "1.205"
$ws.Range("E51").Value = "  +0.46%  "

# Restore the default style on column D so only the underlying value changed
# (the original workbook used the default/unstyled format for these cells).
$ws.Range("D2:D51").Style = "Normal"
